# Applies the numeric updates to the per-job profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as captured by the scheduled market-data refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1588.9286
$ws.Range("I33").Value = 2105.5
$ws.Range("J33").Value = 900.1667
$ws.Range("K33").Value = 2105.5
$ws.Range("L33").Value = 900.1667
$ws.Range("M33").Value = -1876.5
$ws.Range("N33").Value = -1358.1667
$ws.Range("H34").Value = 10397
$ws.Range("I34").Value = 10397
$ws.Range("K34").Value = 10397
$ws.Range("M34").Value = -10194
$ws.Range("H36").Value = 10397
$ws.Range("I36").Value = 10397
$ws.Range("K36").Value = 10397
$ws.Range("M36").Value = -9682
$ws.Range("H70").Value = 36460836
$ws.Range("J70").Value = 25643380
$ws.Range("L70").Value = 76930140
$ws.Range("N70").Value = -76930680
$ws.Range("H73").Value = 36460836
$ws.Range("J73").Value = 25643380
$ws.Range("L73").Value = 76930140
$ws.Range("N73").Value = -76932012
$ws.Range("H87").Value = 76666.336
$ws.Range("J87").Value = 76666.336
$ws.Range("L87").Value = 76666.336
$ws.Range("N87").Value = -79162.336
$ws.Range("H90").Value = 76666.336
$ws.Range("J90").Value = 76666.336
$ws.Range("L90").Value = 229999.008
$ws.Range("N90").Value = -242479.008
$ws.Range("H96").Value = 1164.8334
$ws.Range("I96").Value = 998.75
$ws.Range("J96").Value = 1497
$ws.Range("K96").Value = 2996.25
$ws.Range("L96").Value = 4491
$ws.Range("M96").Value = -1623.25
$ws.Range("N96").Value = -7237
$ws.Range("H98").Value = 3415.8667
$ws.Range("I98").Value = 3415.8667
$ws.Range("K98").Value = 3415.8667
$ws.Range("M98").Value = -1917.8667
$ws.Range("H100").Value = 2996.5
$ws.Range("I100").Value = 1749.5
$ws.Range("J100").Value = 3827.8333
$ws.Range("K100").Value = 1749.5
$ws.Range("L100").Value = 3827.8333
$ws.Range("M100").Value = -1208.5
$ws.Range("N100").Value = -4909.8333
$ws.Range("H103").Value = 1019.17645
$ws.Range("J103").Value = 1129.5
$ws.Range("L103").Value = 3388.5
$ws.Range("N103").Value = -4560.5
$ws.Range("H113").Value = 115390440
$ws.Range("I113").Value = 166668640
$ws.Range("J113").Value = 71437704
$ws.Range("K113").Value = 166668640
$ws.Range("L113").Value = 71437704
$ws.Range("M113").Value = -166665386
$ws.Range("N113").Value = -71444212
$ws.Range("H122").Value = 3415.8667
$ws.Range("I122").Value = 3415.8667
$ws.Range("K122").Value = 10247.6001
$ws.Range("M122").Value = -7797.6001
$ws.Range("H132").Value = 1133.54
$ws.Range("I132").Value = 1115.8776
$ws.Range("K132").Value = 3347.6328
$ws.Range("M132").Value = -817.6328000000003
$ws.Range("H135").Value = 400693.4
$ws.Range("I135").Value = 435449.34
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 3919044.06
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -3916509.06
$ws.Range("N135").Value = -14070
$ws.Range("H137").Value = 4225.0625
$ws.Range("I137").Value = 5600.1665
$ws.Range("K137").Value = 16800.4995
$ws.Range("M137").Value = -14250.4995
$ws.Range("H138").Value = 1767
$ws.Range("J138").Value = 3166.6667
$ws.Range("L138").Value = 9500.000100000001
$ws.Range("N138").Value = -19780.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 50915.883
$ws.Range("I74").Value = 91294.5
$ws.Range("K74").Value = 91294.5
$ws.Range("M74").Value = -90420.5
$ws.Range("H77").Value = 50915.883
$ws.Range("I77").Value = 91294.5
$ws.Range("K77").Value = 456472.5
$ws.Range("M77").Value = -452104.5
$ws.Range("H122").Value = 16999.8
$ws.Range("I122").Value = 23833
$ws.Range("J122").Value = 6750
$ws.Range("K122").Value = 71499
$ws.Range("L122").Value = 20250
$ws.Range("M122").Value = -69049
$ws.Range("N122").Value = -25150
$ws.Range("H132").Value = 7854.4375
$ws.Range("I132").Value = 6404.1904
$ws.Range("K132").Value = 19212.5712
$ws.Range("M132").Value = -16682.5712

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 91005040
$ws.Range("I86").Value = 203680
$ws.Range("J86").Value = 166672830
$ws.Range("K86").Value = 203680
$ws.Range("L86").Value = 166672830
$ws.Range("M86").Value = -202557
$ws.Range("N86").Value = -166675076
$ws.Range("H89").Value = 91005040
$ws.Range("I89").Value = 203680
$ws.Range("J89").Value = 166672830
$ws.Range("K89").Value = 1018400
$ws.Range("L89").Value = 833364150
$ws.Range("M89").Value = -1012784
$ws.Range("N89").Value = -833375382
$ws.Range("H94").Value = 1447.4762
$ws.Range("I94").Value = 614.6667
$ws.Range("J94").Value = 2557.889
$ws.Range("K94").Value = 614.6667
$ws.Range("L94").Value = 2557.889
$ws.Range("M94").Value = -163.6667
$ws.Range("N94").Value = -3459.889
$ws.Range("H134").Value = 5067.9814
$ws.Range("I134").Value = 2022.2424
$ws.Range("K134").Value = 6066.7272
$ws.Range("M134").Value = -3531.7272

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 5000
$ws.Range("J8").Value = 5000
$ws.Range("L8").Value = 5000
$ws.Range("N8").Value = -5280
$ws.Range("H14").Value = 2000
$ws.Range("J14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2340
$ws.Range("H16").Value = 2661.0303
$ws.Range("I16").Value = 1219.409
$ws.Range("K16").Value = 1219.409
$ws.Range("M16").Value = -932.4090000000001
$ws.Range("H31").Value = 9017042
$ws.Range("I31").Value = 3398.0588
$ws.Range("K31").Value = 3398.0588
$ws.Range("M31").Value = -3103.0588
$ws.Range("H34").Value = 9017042
$ws.Range("I34").Value = 3398.0588
$ws.Range("K34").Value = 3398.0588
$ws.Range("M34").Value = -3196.0588
$ws.Range("H58").Value = 12505972
$ws.Range("I58").Value = 27779200
$ws.Range("J58").Value = 9695.5
$ws.Range("K58").Value = 27779200
$ws.Range("L58").Value = 9695.5
$ws.Range("M58").Value = -27778997
$ws.Range("N58").Value = -10101.5
$ws.Range("H105").Value = 10209724
$ws.Range("I105").Value = 17858766
$ws.Range("K105").Value = 17858766
$ws.Range("M105").Value = -17857019
$ws.Range("H113").Value = 2661.0303
$ws.Range("I113").Value = 1219.409
$ws.Range("K113").Value = 1219.409
$ws.Range("M113").Value = 950.5909999999999
$ws.Range("H136").Value = 12505972
$ws.Range("I136").Value = 27779200
$ws.Range("J136").Value = 9695.5
$ws.Range("K136").Value = 83337600
$ws.Range("L136").Value = 29086.5
$ws.Range("M136").Value = -83335050
$ws.Range("N136").Value = -34186.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 47009330
$ws.Range("I109").Value = 45455404
$ws.Range("J109").Value = 55555904
$ws.Range("K109").Value = 136366212
$ws.Range("L109").Value = 166667712
$ws.Range("M109").Value = -136365172
$ws.Range("N109").Value = -166669792
$ws.Range("H124").Value = 3923
$ws.Range("I124").Value = 3923
$ws.Range("K124").Value = 11769
$ws.Range("M124").Value = -6859

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 13419.75
$ws.Range("I36").Value = 13419.75
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 13419.75
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -12934.75
$ws.Range("N36").ClearContents()
$ws.Range("H122").Value = 2417198
$ws.Range("I122").Value = 3623574
$ws.Range("K122").Value = 10870722
$ws.Range("M122").Value = -10868272
$ws.Range("H132").Value = 6284.143
$ws.Range("I132").Value = 2033.7142
$ws.Range("J132").Value = 10534.571
$ws.Range("K132").Value = 6101.142599999999
$ws.Range("L132").Value = 31603.713
$ws.Range("M132").Value = -3571.142599999999
$ws.Range("N132").Value = -36663.713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H68").Value = 3164.5715
$ws.Range("I68").Value = 2033.4445
$ws.Range("K68").Value = 2033.4445
$ws.Range("M68").Value = -1284.4445
$ws.Range("H71").Value = 3164.5715
$ws.Range("I71").Value = 2033.4445
$ws.Range("K71").Value = 10167.2225
$ws.Range("M71").Value = -6423.2225
$ws.Range("H122").Value = 4513
$ws.Range("I122").Value = 3179.8462
$ws.Range("J122").Value = 5957.25
$ws.Range("K122").Value = 9539.5386
$ws.Range("L122").Value = 17871.75
$ws.Range("M122").Value = -7089.5386
$ws.Range("N122").Value = -22771.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 4000
$ws.Range("I6").Value = 4000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 4000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -3885
$ws.Range("N6").ClearContents()
$ws.Range("H9").Value = 10000
$ws.Range("I9").Value = 10000
$ws.Range("K9").Value = 10000
$ws.Range("M9").Value = -9860
$ws.Range("H96").Value = 2428.5715
$ws.Range("I96").Value = 3222
$ws.Range("J96").Value = 1370.6666
$ws.Range("K96").Value = 3222
$ws.Range("L96").Value = 1370.6666
$ws.Range("M96").Value = -1849
$ws.Range("N96").Value = -4116.6666
$ws.Range("H100").Value = 1000.34784
$ws.Range("I100").Value = 833.0909
$ws.Range("J100").Value = 1153.6666
$ws.Range("K100").Value = 1666.1818
$ws.Range("L100").Value = 2307.3332
$ws.Range("M100").Value = -1125.1818
$ws.Range("N100").Value = -3389.3332
$ws.Range("H107").Value = 555.4583
$ws.Range("I107").Value = 381.35
$ws.Range("J107").Value = 1426
$ws.Range("K107").Value = 1144.05
$ws.Range("L107").Value = 4278
$ws.Range("M107").Value = 775.9499999999998
$ws.Range("N107").Value = -8118
$ws.Range("H136").Value = 22754486
$ws.Range("I136").Value = 47620068
$ws.Range("K136").Value = 142860204
$ws.Range("M136").Value = -142857654
